# Commit: "update coupon code for MWL"
#
# The source workbook keeps a pool of previously-generated coupon codes in
# xl/sharedStrings.xml; the active coupon for this row (Sheet1!A2) is swapped
# to a freshly generated code. Apply the visible, COM-addressable part of
# that change: update the coupon code cell itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CA-VZM6NS56"
